# Adds the Artisan Command "adjustSV" to the Commands sheet of eventbuttons.xlsx
# and updates the existing "pidSV" command to use an <int> argument instead of
# a <float> one (to match the new adjustSV(<int>) sibling command).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert a brand new row right above the existing "pidSV(<float>)" row (row 79)
# so that all rows below shift down by one.
$ws.Rows.Item(79).Insert()

# New row 79: the adjustSV command and its description.
$ws.Range("B79").Value = "adjustSV(<int>)"
$ws.Range("C79").Value = "increases or decreases the current target SV value by <int>"
$ws.Rows.Item(79).RowHeight = 13.8

# The old pidSV row now lives at row 80 (its description in C80 shifted down
# automatically); update its command signature from <float> to <int>.
$ws.Range("B80").Value = "pidSV(<int>)"

Write-Output "adjustSV command added; pidSV argument updated to <int>"
